# Delete the "Homepage-0" worksheet entirely
$wb = $excel.ActiveWorkbook

$wsHome0 = $wb.Worksheets.Item("Homepage-0")
$wsHome0.Delete()

# On "Homepage-6", remove the 3 duplicate rows that hold the
# "happy-mothers-day-from-johnson-johnson" URLs (rows 125-127)
$wsHome6 = $wb.Worksheets.Item("Homepage-6")
$rng = $wsHome6.Range("A125:A127")
$rng.EntireRow.Delete()

# Autofit column A to match the "bestFit" width Excel would apply
$wsHome6.Columns.Item(1).AutoFit()

# Make Homepage-6 the active sheet (matches activeTab state in the diff)
$wsHome6.Activate()
$wsHome6.Range("A125").Select()

$wb.Save()
